# Update attendance ("想去人数") figures in 展览 and 全部类型 sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13088
$ws1.Range("F10").Value = 13045
$ws1.Range("F11").Value = 301
$ws1.Range("F14").Value = 7774
$ws1.Range("F15").Value = 213
$ws1.Range("F16").Value = 132
$ws1.Range("F26").Value = 5219

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13088
$ws4.Range("F11").Value = 13045
$ws4.Range("F12").Value = 301
$ws4.Range("F15").Value = 7774
$ws4.Range("F16").Value = 213
$ws4.Range("F17").Value = 132
$ws4.Range("F29").Value = 5219
